# Atualização automática via cronjob
# Refresh the "vendas atipicas" dataset: replace the data rows (2..21) with a
# new, smaller dataset (2..15) reflecting the latest cronjob run, and shrink
# the used range accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1,  "2025-03-26", 40,  "V V REFEICOES LTDA",                               "000091", "VASSOURA VARRE CANTO COM CABO PLASTIFICADO",                   65,   $false),
    @(3,  "2025-03-26", 10,  "JURUA ESTALEIROS E NAVEGACAO LTDA",                "000425", "COADOR DE CAFE EG (EXTRA GRANDE)",                             8,    $false),
    @(8,  "2025-03-26", 200, "JURUA ESTALEIROS E NAVEGACAO LTDA",                "000122", "SABAO EM PO ALA LAVANDA ROUPAS 400G",                          296,  $false),
    @(12, "2025-03-28", 350, "AMAZONPEL COMERCIO DE MATERIAIS DE LIMPEZA LTDA",  "000494", "FIBRA DE LIMPEZA PESADA 98X229MM SLIM NOBRE",                  331,  $true),
    @(6,  "2025-04-01", 20,  "SAT BRAS INDUSTRIA ELETRONICA DA AMAZONIA LTDA.",  "000258", "DISPENSER PAPEL HIGIENICO ROLAO 300-500M BRANCO NOBRE STREET", 40,   $true),
    @(11, "2025-04-01", 61,  "AMAZONIA FORMULA LTDA",                            "000799", "SACO DE LIXO 30L REFORCADO PACOTINHO C/10 UND",                352,  $false),
    @(0,  "2025-04-02", 50,  "REVEMAR COMERCIO DE MOTOS LTDA",                   "000054", "PAPEL HIGIENICO 8X300 NEWPAPER 100% Celulose",                 206,  $false),
    @(2,  "2025-04-02", 250, "RH MULTI SERVICOS ADMINISTRATIVOS S.A",            "000041", "LUVAS DESCARTAVEIS C/ 100 UND",                                1632, $false),
    @(4,  "2025-04-02", 100, "REVEMAR COMERCIO DE MOTOS LTDA",                   "010041", "PAPEL A4 REPORT RESMA C 500FLS",                               41,   $false),
    @(5,  "2025-04-02", 30,  "REVEMAR COMERCIO DE MOTOS LTDA",                   "000434", "FRASCO COM VALVULA PUMP 450ML (p/alcool/sabonete) - NOBRE",    35,   $true),
    @(7,  "2025-04-02", 150, "REVEMAR COMERCIO DE MOTOS LTDA",                   "000349", "DESODORISADOR LADY AEROSSOL 360ML TALCO SUAVE CARINHO",        820,  $true),
    @(9,  "2025-04-02", 100, "REVEMAR COMERCIO DE MOTOS LTDA",                   "000057", "PAPEL TOALHA 8x100M NEWPAPER 100% CELULOSE",                   75,   $false),
    @(13, "2025-04-02", 15,  "REVEMAR COMERCIO DE MOTOS LTDA",                   "010180", "GRAMPEADOR METAL 25FLS 11,5 CM JOCAR OFFICE",                  1,    $false),
    @(10, "2025-04-03", 120, "SIND.DAS EMPRESAS DE TRANSP.DE PASSAGEIROS DO EST.","000782","ALCOOL LIQUIDO 70% INPM SANTA CRUZ 1L",                        10,   $false)
)

$newLastRow = 1 + $data.Count   # header is row 1
$oldLastRow = $ws.UsedRange.Rows.Count
$xlPasteValues = -4163

# Write the refreshed data into rows 2..($newLastRow)
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]

    # "Dia" (B) holds date-like text ("2025-03-26") and "id_produto" (E)
    # holds zero-padded numeric codes ("000091"); a plain .Value assignment
    # would make Excel reinterpret these as a real date serial / number,
    # losing the original text representation. Route the literal text
    # through a throwaway formula and paste back only the resulting value,
    # which keeps the cell a plain shared-string without touching styles.
    $ws.Cells.Item($r, 2).Formula = '="' + $row[1] + '"'
    $ws.Cells.Item($r, 2).Copy()
    $ws.Cells.Item($r, 2).PasteSpecial($xlPasteValues)

    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]

    $ws.Cells.Item($r, 5).Formula = '="' + $row[4] + '"'
    $ws.Cells.Item($r, 5).Copy()
    $ws.Cells.Item($r, 5).PasteSpecial($xlPasteValues)

    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}
$excel.CutCopyMode = $false

# Remove the now-obsolete trailing rows so the sheet shrinks to A1:H15
if ($oldLastRow -gt $newLastRow) {
    $deleteRange = $ws.Range($ws.Cells.Item($newLastRow + 1, 1), $ws.Cells.Item($oldLastRow, 8))
    $deleteRange.Delete()
}

$wb.Save()
